# Apply the "fix command names, clarify startup state for prepackaged VM" edit
$d = $word.ActiveDocument

# --- Resize the "Heads up" text box first (layout autofit delta picked up by
#     later edits would otherwise clobber this), matching the new measured
#     extent recorded in the target document. ---
$s = $d.Shapes.Item(1)
$s.Width = 430.4409448818898
$s.Height = 37.75748031496063

# --- Step 1: reword the "Boot your Linux system or VM..." instruction ---
$d.Content.Find.Execute("Boot your Linux system or VM.  If needed, login and open a terminal and cd to the " + [char]34 + "~/labtainer/labtainer-student" + [char]34 + " directory.  (The prepackaged Labtainer VM starts with such a terminal.)  Then start the lab:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Boot your Linux system or VM.  If necessary, log in and then open a terminal window and cd to the labtainer/labtainer-student directory.  The pre-packaged Labtainer VM will start with such a terminal open for you.   Then start the lab:", 2)

# --- Step 2: rename the start command ---
$d.Content.Find.Execute("start.py centos-log", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "labtainer centos-log", 2)

# --- Step 3: re-apply the "Note the original terminal..." sentence so the
#     three runs coalesce into one (content is unchanged) ---
$d.Content.Find.Execute("Note the original terminal displays the paths to two files on your Linux host: ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Note the original terminal displays the paths to two files on your Linux host: ", 2)

# --- Step 4: stop.py -> stoplab in the "copying the completed report back" sentence ---
$d.Content.Find.Execute([char]0x201c + "stop.py" + [char]0x201d + " to stop the lab for the last time.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    [char]0x201c + "stoplab" + [char]0x201d + " to stop the lab for the last time.", 2)

# --- Step 5: standalone "stop.py " command line -> "stoplab " ---
$d.Content.Find.Execute("stop.py ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "stoplab ", 2)

# --- Step 6: "./stop.py" -> "stoplab" in the "must copy that completed file" sentence ---
$d.Content.Find.Execute([char]0x201c + "./stop.py" + [char]0x201d + ".", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    [char]0x201c + "stoplab" + [char]0x201d + ".", 2)

# --- Step 7: widen the appendix table's left cell margin (133 -> 143 dxa) ---
$t = $d.Tables.Item(1)
$t.LeftPadding = 7.15
